$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T220004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2006%3A00%3A04%20GMT&X-Amz-Signature=111ba03a86cfc68647a8f16801c3d046e4066c0510041d2d68d1f2e6151edbcd"
$ws.Range("K3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T220004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2006%3A00%3A04%20GMT&X-Amz-Signature=c02eb2a5b6b8464de9faa93e5a258a3dc1cb58e7b6a187804cbca891518d40f2"
$ws.Range("K4").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T220004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2006%3A00%3A04%20GMT&X-Amz-Signature=8173747d578165e993f508c27c486c38d432e1aeedb1ffdf99086a25bedf3a16"
$ws.Range("K5").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T220004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2006%3A00%3A04%20GMT&X-Amz-Signature=48e9c22f9873ca6d1f8fd8599610ad3610afe6d1a967a85755635932bed78e5a"
$ws.Range("K6").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T220004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2006%3A00%3A04%20GMT&X-Amz-Signature=30d3c301af36090453bcf397fe36d91d877f9be3d9119a4132c4149baf930680"

